$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 16500
$ws.Range("F2").Value = 66500

# Row 3
$ws.Range("B3").Value = 13500
$ws.Range("C3").Value = 20000
$ws.Range("D3").Value = 20000
$ws.Range("E3").Value = 10000
$ws.Range("F3").Value = 63500

# Row 4
$ws.Range("D4").Value = 20000
$ws.Range("F4").Value = 42500

# Row 5
$ws.Range("D5").Value = 20000
$ws.Range("F5").Value = 42500

# Row 7
$ws.Range("D7").Value = 0
$ws.Range("F7").Value = 22500

# Row 8
$ws.Range("D8").Value = 0
$ws.Range("F8").Value = 22500

# Row 9
$ws.Range("D9").Value = 0
$ws.Range("F9").Value = 22500

# Row 10
$ws.Range("D10").Value = 0
$ws.Range("F10").Value = 22500

# Row 11
$ws.Range("D11").Value = 10000
$ws.Range("F11").Value = 32500

# Row 12
$ws.Range("D12").Value = 10000
$ws.Range("F12").Value = 32500

# Row 13
$ws.Range("D13").Value = 10000
$ws.Range("F13").Value = 32500

# Row 14
$ws.Range("D14").Value = 10000
$ws.Range("F14").Value = 32500

# Row 15
$ws.Range("B15").Value = 12500
$ws.Range("C15").Value = 10000
$ws.Range("F15").Value = 52500

# Row 24
$ws.Range("B24").Value = 22000
$ws.Range("F24").Value = 72000

# Row 25
$ws.Range("B25").Value = 12500
$ws.Range("F25").Value = 52500
